# Applies the "Corrected excel sheets for application fix issues" edit.
$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "Summary"
# ---------------------------------------------------------------------
$wsSummary = $wb.Worksheets.Item("Summary")
$wsSummary.Range("F2").Value = 0
$wsSummary.Range("A3").Value = 720.4
$wsSummary.Range("E3").Value = 520.4

# ---------------------------------------------------------------------
# Sheet "Repayment schedule"
# ---------------------------------------------------------------------
$wsSchedule = $wb.Worksheets.Item("Repayment schedule")

# Row 2 - drop the stray P2 cell entirely (no longer populated)
$wsSchedule.Range("P2").Clear()

# Row 3 - drop O3 entirely (column O no longer populated below the header row)
$wsSchedule.Range("O3").Clear()

# Row 4
$wsSchedule.Range("B4").Value = 31
$wsSchedule.Range("C4").Value = 42095
$wsSchedule.Range("F4").Value = 872.19
$wsSchedule.Range("G4").Value = 8363.27
$wsSchedule.Range("H4").Value = 92.35
$wsSchedule.Range("O4").Clear()

# Row 5
$wsSchedule.Range("B5").Value = 30
$wsSchedule.Range("C5").Value = 42125
$wsSchedule.Range("F5").Value = 880.91
$wsSchedule.Range("G5").Value = 7482.36
$wsSchedule.Range("H5").Value = 83.63
$wsSchedule.Range("O5").Clear()

# Row 6
$wsSchedule.Range("B6").Value = 31
$wsSchedule.Range("C6").Value = 42156
$wsSchedule.Range("F6").Value = 889.72
$wsSchedule.Range("G6").Value = 6592.64
$wsSchedule.Range("H6").Value = 74.82
$wsSchedule.Range("O6").Clear()

# Row 7
$wsSchedule.Range("B7").Value = 30
$wsSchedule.Range("C7").Value = 42186
$wsSchedule.Range("F7").Value = 898.61
$wsSchedule.Range("G7").Value = 5694.03
$wsSchedule.Range("H7").Value = 65.93
$wsSchedule.Range("O7").Clear()

# Row 8
$wsSchedule.Range("B8").Value = 31
$wsSchedule.Range("C8").Value = 42217
$wsSchedule.Range("F8").Value = 907.6
$wsSchedule.Range("G8").Value = 4786.43
$wsSchedule.Range("H8").Value = 56.94
$wsSchedule.Range("O8").Clear()

# Row 9
$wsSchedule.Range("C9").Value = 42248
$wsSchedule.Range("F9").Value = 916.68
$wsSchedule.Range("G9").Value = 3869.75
$wsSchedule.Range("H9").Value = 47.86
$wsSchedule.Range("O9").Clear()

# Row 10
$wsSchedule.Range("B10").Value = 30
$wsSchedule.Range("C10").Value = 42278
$wsSchedule.Range("F10").Value = 925.84
$wsSchedule.Range("G10").Value = 2943.91
$wsSchedule.Range("H10").Value = 38.7
$wsSchedule.Range("O10").Clear()

# Row 11
$wsSchedule.Range("B11").Value = 31
$wsSchedule.Range("C11").Value = 42309
$wsSchedule.Range("F11").Value = 935.1
$wsSchedule.Range("G11").Value = 2008.81
$wsSchedule.Range("H11").Value = 29.44
$wsSchedule.Range("O11").Clear()

# Row 12
$wsSchedule.Range("B12").Value = 30
$wsSchedule.Range("C12").Value = 42339
$wsSchedule.Range("F12").Value = 944.45
$wsSchedule.Range("G12").Value = 1064.36
$wsSchedule.Range("G12").NumberFormat = "#,##0.00"
$wsSchedule.Range("H12").Value = 20.09
$wsSchedule.Range("O12").Clear()

# Row 13
$wsSchedule.Range("B13").Value = 31
$wsSchedule.Range("C13").Value = 42370
$wsSchedule.Range("F13").Value = 1064.36
$wsSchedule.Range("F13").NumberFormat = "#,##0.00"
$wsSchedule.Range("H13").Value = 10.64
$wsSchedule.Range("K13").Value = 1075
$wsSchedule.Range("K13").NumberFormat = "#,##0"
$wsSchedule.Range("O13").Clear()
$wsSchedule.Range("P13").Value = 1075
$wsSchedule.Range("P13").NumberFormat = "#,##0"

# ---------------------------------------------------------------------
# Sheet "Transactions"
# ---------------------------------------------------------------------
$wsTransactions = $wb.Worksheets.Item("Transactions")
$wsTransactions.Range("A2").Value = 6352
$wsTransactions.Range("A3").Value = 6350

# ---------------------------------------------------------------------
# Restore per-sheet cursor/selection to match the saved view state
# ---------------------------------------------------------------------
$wsSummary.Activate()
$wsSummary.Range("C4").Select()

$wsSchedule.Activate()
$wsSchedule.Range("F4:F13").Select()

$wsTransactions.Activate()
$wsTransactions.Range("D3").Select()
